# "Join en Create menu aanmaken + button working"
#
# Adds two new fields ("naam" and "code") to the small "winnaar / diamant /
# positie in ranglijst" table in column C (rows 7-12) of the single sheet,
# and leaves the selection on C10 (where the user ended up after typing the
# new values in).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared-string values, mirroring the existing "winnaar" (C7) entry
# immediately above them in the same mini-table.
$ws.Range("C8").Value = "naam"
$ws.Range("C9").Value = "code"

# Match the (non-bold) look of the sibling cell C7 directly above this pair,
# so the two new rows read the same as the rest of the little table.
$ws.Range("C8:C9").Font.Bold = $false

# Final cursor position left behind by the edit.
$ws.Range("C10").Select()
